# Auto-generated edit script: updates market-board derived profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled pricing refresh.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 962.625
$ws.Range("J32").Value = 1425
$ws.Range("L32").Value = 1425
$ws.Range("N32").Value = -2077
$ws.Range("H47").Value = 47537
$ws.Range("J47").Value = 47537
$ws.Range("L47").Value = 47537
$ws.Range("N47").Value = -49481
$ws.Range("H95").Value = 16000
$ws.Range("J95").Value = 16000
$ws.Range("L95").Value = 16000
$ws.Range("N95").Value = -21492
$ws.Range("H138").Value = 2117.4424
$ws.Range("I138").Value = 899.93335
$ws.Range("J138").Value = 2611.027
$ws.Range("K138").Value = 2699.80005
$ws.Range("L138").Value = 7833.081
$ws.Range("M138").Value = 2440.19995
$ws.Range("N138").Value = -18113.081

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 252506.9
$ws.Range("I32").Value = 283441.2
$ws.Range("J32").Value = 128769.78
$ws.Range("K32").Value = 283441.2
$ws.Range("L32").Value = 128769.78
$ws.Range("M32").Value = -283154.2
$ws.Range("N32").Value = -129343.78
$ws.Range("H61").Value = 1810.8372
$ws.Range("I61").Value = 1734.8572
$ws.Range("K61").Value = 1734.8572
$ws.Range("M61").Value = -1522.8572
$ws.Range("H74").Value = 193815.25
$ws.Range("I74").Value = 358710.4
$ws.Range("J74").Value = 1437.5834
$ws.Range("K74").Value = 358710.4
$ws.Range("L74").Value = 1437.5834
$ws.Range("M74").Value = -357836.4
$ws.Range("N74").Value = -3185.5834
$ws.Range("H77").Value = 193815.25
$ws.Range("I77").Value = 358710.4
$ws.Range("J77").Value = 1437.5834
$ws.Range("K77").Value = 1793552
$ws.Range("L77").Value = 7187.916999999999
$ws.Range("M77").Value = -1789184
$ws.Range("N77").Value = -15923.917
$ws.Range("H95").Value = 26500
$ws.Range("J95").Value = 26500
$ws.Range("L95").Value = 26500
$ws.Range("N95").Value = -31992
$ws.Range("H110").Value = 788.7857
$ws.Range("I110").Value = 788.7857
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 788.7857
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1256.2143
$ws.Range("N110").ClearContents()
$ws.Range("H136").Value = 1810.8372
$ws.Range("I136").Value = 1734.8572
$ws.Range("K136").Value = 5204.571599999999
$ws.Range("M136").Value = -2654.571599999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1183.9231
$ws.Range("I107").Value = 1032.5834
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 1032.5834
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 887.4166
$ws.Range("N107").Value = -6840
$ws.Range("H134").Value = 12507563
$ws.Range("I134").Value = 16135984
$ws.Range("J134").Value = 9668.223
$ws.Range("K134").Value = 48407952
$ws.Range("L134").Value = 29004.669
$ws.Range("M134").Value = -48405417
$ws.Range("N134").Value = -34074.669

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2682
$ws.Range("I22").Value = 2682
$ws.Range("K22").Value = 2682
$ws.Range("M22").Value = -2332
$ws.Range("H31").Value = 5112.5
$ws.Range("I31").Value = 3909.0303
$ws.Range("J31").Value = 6701.08
$ws.Range("K31").Value = 3909.0303
$ws.Range("L31").Value = 6701.08
$ws.Range("M31").Value = -3614.0303
$ws.Range("N31").Value = -7291.08
$ws.Range("H34").Value = 5112.5
$ws.Range("I34").Value = 3909.0303
$ws.Range("J34").Value = 6701.08
$ws.Range("K34").Value = 3909.0303
$ws.Range("L34").Value = 6701.08
$ws.Range("M34").Value = -3707.0303
$ws.Range("N34").Value = -7105.08
$ws.Range("H58").Value = 3531.4285
$ws.Range("I58").Value = 1167
$ws.Range("K58").Value = 1167
$ws.Range("M58").Value = -964
$ws.Range("H136").Value = 3531.4285
$ws.Range("I136").Value = 1167
$ws.Range("K136").Value = 3501
$ws.Range("M136").Value = -951

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1455
$ws.Range("I22").Value = 790
$ws.Range("J22").Value = 1588
$ws.Range("K22").Value = 2370
$ws.Range("L22").Value = 4764
$ws.Range("M22").Value = -2201
$ws.Range("N22").Value = -5102
$ws.Range("H27").Value = 1455
$ws.Range("I27").Value = 790
$ws.Range("J27").Value = 1588
$ws.Range("K27").Value = 2370
$ws.Range("L27").Value = 4764
$ws.Range("M27").Value = -2268
$ws.Range("N27").Value = -4968
$ws.Range("H68").Value = 952.5282999999999
$ws.Range("J68").Value = 1045.9025
$ws.Range("L68").Value = 3137.7075
$ws.Range("N68").Value = -4759.7075
$ws.Range("H71").Value = 952.5282999999999
$ws.Range("J71").Value = 1045.9025
$ws.Range("L71").Value = 9413.122499999999
$ws.Range("N71").Value = -17525.1225
$ws.Range("H86").Value = 265
$ws.Range("I86").Value = 231.25
$ws.Range("J86").Value = 400
$ws.Range("K86").Value = 693.75
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = 492.25
$ws.Range("N86").Value = -3572
$ws.Range("H89").Value = 265
$ws.Range("I89").Value = 231.25
$ws.Range("J89").Value = 400
$ws.Range("K89").Value = 2081.25
$ws.Range("L89").Value = 3600
$ws.Range("M89").Value = 3846.75
$ws.Range("N89").Value = -15456
$ws.Range("H107").Value = 454.01755
$ws.Range("I107").Value = 240.3158
$ws.Range("J107").Value = 881.4211
$ws.Range("K107").Value = 720.9474
$ws.Range("L107").Value = 2644.2633
$ws.Range("M107").Value = 1199.0526
$ws.Range("N107").Value = -6484.263300000001
$ws.Range("H131").Value = 954.1539
$ws.Range("I131").Value = 441.5
$ws.Range("J131").Value = 1026.1052
$ws.Range("K131").Value = 1324.5
$ws.Range("L131").Value = 3078.3156
$ws.Range("M131").Value = 3715.5
$ws.Range("N131").Value = -13158.3156

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 18818
$ws.Range("J15").Value = 18818
$ws.Range("L15").Value = 18818
$ws.Range("N15").Value = -19394
$ws.Range("H81").Value = 18818
$ws.Range("J81").Value = 18818
$ws.Range("L81").Value = 18818
$ws.Range("N81").Value = -20814
$ws.Range("H84").Value = 18818
$ws.Range("J84").Value = 18818
$ws.Range("L84").Value = 56454
$ws.Range("N84").Value = -66438
$ws.Range("H122").Value = 2901
$ws.Range("I122").Value = 2670.6
$ws.Range("J122").Value = 3477
$ws.Range("K122").Value = 8011.799999999999
$ws.Range("L122").Value = 10431
$ws.Range("M122").Value = -5561.799999999999
$ws.Range("N122").Value = -15331
$ws.Range("H132").Value = 22826.082
$ws.Range("I132").Value = 1934.9697
$ws.Range("J132").Value = 65914
$ws.Range("K132").Value = 5804.909100000001
$ws.Range("L132").Value = 197742
$ws.Range("M132").Value = -3274.909100000001
$ws.Range("N132").Value = -202802

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1629.1428
$ws.Range("I22").Value = 2001
$ws.Range("J22").Value = 1567.1666
$ws.Range("K22").Value = 2001
$ws.Range("L22").Value = 1567.1666
$ws.Range("M22").Value = -1706
$ws.Range("N22").Value = -2157.1666
$ws.Range("H27").Value = 1629.1428
$ws.Range("I27").Value = 2001
$ws.Range("J27").Value = 1567.1666
$ws.Range("K27").Value = 2001
$ws.Range("L27").Value = 1567.1666
$ws.Range("M27").Value = -1894
$ws.Range("N27").Value = -1781.1666
$ws.Range("H68").Value = 2076.25
$ws.Range("I68").Value = 1775
$ws.Range("K68").Value = 1775
$ws.Range("M68").Value = -1026
$ws.Range("H71").Value = 2076.25
$ws.Range("I71").Value = 1775
$ws.Range("K71").Value = 8875
$ws.Range("M71").Value = -5131

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 12501.5
$ws.Range("J2").Value = 12501.5
$ws.Range("L2").Value = 12501.5
$ws.Range("N2").Value = -12725.5
$ws.Range("H81").Value = 2464.2856
$ws.Range("I81").Value = 860
$ws.Range("K81").Value = 1720
$ws.Range("M81").Value = -659
$ws.Range("H84").Value = 2464.2856
$ws.Range("I84").Value = 860
$ws.Range("K84").Value = 8600
$ws.Range("M84").Value = -3296
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H136").Value = 287134.47
$ws.Range("I136").Value = 1860.1052
$ws.Range("J136").Value = 625897.75
$ws.Range("K136").Value = 5580.3156
$ws.Range("L136").Value = 1877693.25
$ws.Range("M136").Value = -3030.3156
$ws.Range("N136").Value = -1882793.25
